# Project Sample Project is saved. Author: admin. Type: SAVE.
# Change the "Return" value of the first decision-table rule row (E8)
# from "Good Morning" to "Good Morning2".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("E8").Value = "Good Morning2"
